# fix: Change the position of MEMBERS in template.xlsx
#
# The two reminder-message templates on the "config" sheet had "{MEMBERS}"
# at the very start of the text. Move it to the end instead (after
# "{ZOOM_TEXT}"), and drop the now-redundant blank line that used to sit
# between "{SCHEDULE}" and "{ZOOM_TEXT}".
#
# Also reflects that the "config" sheet tab became the active one (instead
# of "schedule"), with a new active-cell selection on that sheet.

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("config")

# --- B2: message sent at 9am ---------------------------------------------
$wsConfig.Range("B2").Value = "【リマインド】 本日{SCHEDULE_LEN}つの予定があります。`n{SCHEDULE}`n{ZOOM_TEXT}`n{MEMBERS}"

# --- B3: message sent 1 hour before the event ------------------------------
$wsConfig.Range("B3").Value = "【リマインド】 もうすぐ以下の予定が開始されます。`n{SCHEDULE}`n{ZOOM_TEXT}`n{MEMBERS}"

# Both cells keep their word-wrap formatting; row heights stay as authored.
$wsConfig.Range("B2:B3").WrapText = $true
$wsConfig.Rows.Item(2).RowHeight = 15.75
$wsConfig.Rows.Item(3).RowHeight = 15.75

# --- Active sheet / selection moved from "schedule" to "config" -----------
$wsConfig.Activate() | Out-Null
$wsConfig.Range("C23").Select() | Out-Null
